# "Alterado cor da camada de normalizacao" (Changed color of the normalization layer)
#
# 1) The three "LayerNorm" rectangles on the slide get their fill changed
#    from the theme color accent4 (lumMod 20%/lumOff 80%) to the flat
#    RGB color DEEBF7 (matching the other boxes already using that color).
# 2) The slide-footer date field cached in the slide master + every slide
#    layout is refreshed from 14/11/2023 to 02/05/2024.
# 3) The legend caption "Camadas removidas durante a inferencia" is split
#    into three runs and extended to also mention "+ LayerNorm".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Recolor the three LayerNorm boxes (shape ids 164, 200, 141) to DEEBF7
# ---------------------------------------------------------------------
$layerNormIds = @(164, 200, 141)
$targetRgb = 0xF7EBDE   # PowerPoint RGB() is stored BGR-packed: 0x00DEEBF7 -> R=DE,G=EB,B=F7

$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($layerNormIds -contains $shp.Id) {
        $shp.Fill.ForeColor.RGB = $targetRgb
    }
}

# ---------------------------------------------------------------------
# 2) Refresh the cached "datetimeFigureOut" footer text site-wide
# ---------------------------------------------------------------------
$newDate = "02/05/2024"

function Update-DatePlaceholders($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        $isDate = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) { $isDate = $true }
        } catch {
        }
        if ($isDate -and $sh.HasTextFrame) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

Update-DatePlaceholders $p.SlideMaster.Shapes
for ($L = 1; $L -le $p.SlideMaster.CustomLayouts.Count; $L++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($L)
    Update-DatePlaceholders $layout.Shapes
}

# ---------------------------------------------------------------------
# 3) Extend / split the "Camadas removidas durante a inferencia" caption
# ---------------------------------------------------------------------
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.Id -eq 217) {
        $tr = $shp.TextFrame.TextRange
        $tr.Text = "Camadas removidas durante a inferência + LayerNorm"

        # Re-stamp (no-op) italics on the 2nd/3rd segments so the engine
        # materializes them as independent <a:r> runs, matching the
        # authored run split, while keeping identical run properties.
        $run2 = $tr.Characters(29, 13)   # "inferência + "
        $run2.Font.Italic = $true

        $run3 = $tr.Characters(42, 9)    # "LayerNorm"
        $run3.Font.Italic = $true
    }
}
